# Insert a new data row at row 655 (shifting existing rows 655-713 down to
# 656-714), then populate the new row with its values. This mirrors the
# source diff: dimension grows from A1:T713 to A1:T714, and every row from
# the old 655 onward is pushed down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(655).Insert()

$ws.Cells.Item(655, 1).Value = 9
$ws.Cells.Item(655, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(655, 3).Value = "Metropolitana"
$ws.Cells.Item(655, 4).Value = 45132
$ws.Cells.Item(655, 5).Value = 13
$ws.Cells.Item(655, 6).Value = "Fruta"
$ws.Cells.Item(655, 7).Value = 100108
$ws.Cells.Item(655, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(655, 9).Value = 100108002
$ws.Cells.Item(655, 10).Value = "Mango"
$ws.Cells.Item(655, 11).Value = "Sin especificar"
$ws.Cells.Item(655, 12).Value = "Primera"
$ws.Cells.Item(655, 13).Value = 580
$ws.Cells.Item(655, 14).Value = 7500
$ws.Cells.Item(655, 15).Value = 8000
$ws.Cells.Item(655, 16).Value = 7759
$ws.Cells.Item(655, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(655, 18).Value = "Brasil"
$ws.Cells.Item(655, 19).Value = 1940
$ws.Cells.Item(655, 20).Value = 4
